$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear contents (values) of A3:E6 while preserving cell formatting/styles
$ws.Range("A3:E6").ClearContents()

# Update the active selection to match the authored state (L5)
$ws.Range("L5").Select()
